# Adds COVID-style state data for rows 140-144 (shared strings 194-198, already present)
# and appends 4 new trailing date rows 145-148 with new shared strings 199-202.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 140: fill data columns B:BE (columns E, N, AC, AR, AZ stay blank, as in every other data row)
$row = New-Object 'object[,]' 1,56
$row[0,0] = 3.54933378
$row[0,1] = 22.75420331
$row[0,2] = 29.21816616
$row[0,4] = 32.06064262
$row[0,5] = 7.57006494
$row[0,6] = 6.29778047
$row[0,7] = 11.55155918
$row[0,8] = 14.89273478
$row[0,9] = 23.06657092
$row[0,10] = 14.26889578
$row[0,11] = 13.84721165
$row[0,13] = 0
$row[0,14] = 29.27596068
$row[0,15] = 26.29854249
$row[0,16] = 8.26333337
$row[0,17] = 8.68358276
$row[0,18] = 14.21616393
$row[0,19] = 17.99344383
$row[0,20] = 11.95357201
$row[0,21] = 12.15381626
$row[0,22] = 10.24755617
$row[0,23] = 0
$row[0,24] = 7.36503458
$row[0,25] = 14.97123753
$row[0,26] = 9.24353582
$row[0,28] = 49.82325379
$row[0,29] = 32.32487106
$row[0,30] = 12.03570985
$row[0,31] = 3.70873099
$row[0,32] = 15.57609549
$row[0,33] = 18.33635983
$row[0,34] = 12.0839654
$row[0,35] = 41.07691133
$row[0,36] = 17.8429001
$row[0,37] = 8.573965400000001
$row[0,38] = 8.749599570000001
$row[0,39] = 30.36990321
$row[0,40] = 14.05490297
$row[0,41] = 7.77396204
$row[0,43] = 22.01483012
$row[0,44] = 19.99209226
$row[0,45] = 21.4748459
$row[0,46] = 18.15143556
$row[0,47] = 19.37779344
$row[0,48] = 14.66878021
$row[0,49] = 6.197653
$row[0,51] = 15.38869597
$row[0,52] = 6.65820503
$row[0,53] = 30.66829274
$row[0,54] = 0
$row[0,55] = 0
$ws.Range("B140:BE140").Value = $row

# Row 141: fill data columns B:BE (columns E, N, AC, AR, AZ stay blank, as in every other data row)
$row = New-Object 'object[,]' 1,56
$row[0,0] = 23.07197513
$row[0,1] = 31.51753125
$row[0,2] = 22.81315585
$row[0,4] = 31.31552843
$row[0,5] = 7.85835547
$row[0,6] = 6.14485069
$row[0,7] = 11.12484049
$row[0,8] = 8.074904610000001
$row[0,9] = 19.56387429
$row[0,10] = 12.73009263
$row[0,11] = 17.21042151
$row[0,13] = 0
$row[0,14] = 42.86734835
$row[0,15] = 40.13876664
$row[0,16] = 8.64318411
$row[0,17] = 17.03205038
$row[0,18] = 12.53056235
$row[0,19] = 17.84808536
$row[0,20] = 8.489246659999999
$row[0,21] = 11.9978301
$row[0,22] = 10.69753528
$row[0,23] = 0
$row[0,24] = 6.34390689
$row[0,25] = 17.43779789
$row[0,26] = 15.81642045
$row[0,28] = 40.99815773
$row[0,29] = 26.50919815
$row[0,30] = 13.93130408
$row[0,31] = 2.16727502
$row[0,32] = 12.2306104
$row[0,33] = 13.86215093
$row[0,34] = 12.77141751
$row[0,35] = 47.6737252
$row[0,36] = 22.86265853
$row[0,37] = 8.865683539999999
$row[0,38] = 8.16612054
$row[0,39] = 31.20472873
$row[0,40] = 14.22864228
$row[0,41] = 7.31798952
$row[0,43] = 17.66915876
$row[0,44] = 18.58322389
$row[0,45] = 15.25901598
$row[0,46] = 14.57208722
$row[0,47] = 19.76659741
$row[0,48] = 11.05211921
$row[0,49] = 5.5915531
$row[0,51] = 13.20561135
$row[0,52] = 5.07218123
$row[0,53] = 33.99274262
$row[0,54] = 0
$row[0,55] = 0
$ws.Range("B141:BE141").Value = $row

# Row 142: fill data columns B:BE (columns E, N, AC, AR, AZ stay blank, as in every other data row)
$row = New-Object 'object[,]' 1,56
$row[0,0] = 17.35712836
$row[0,1] = 25.90733783
$row[0,2] = 16.91026042
$row[0,4] = 32.99912113
$row[0,5] = 8.501628849999999
$row[0,6] = 6.48126063
$row[0,7] = 7.69126479
$row[0,8] = 2.6065077
$row[0,9] = 15.91563657
$row[0,10] = 15.01950778
$row[0,11] = 18.31282645
$row[0,13] = 0
$row[0,14] = 46.03606968
$row[0,15] = 33.5671213
$row[0,16] = 8.247448690000001
$row[0,17] = 19.0361658
$row[0,18] = 10.63527391
$row[0,19] = 13.97616171
$row[0,20] = 5.51550835
$row[0,21] = 8.527361170000001
$row[0,22] = 7.46498967
$row[0,23] = 11.73917124
$row[0,24] = 4.1600268
$row[0,25] = 15.38639743
$row[0,26] = 22.91948226
$row[0,28] = 39.37454104
$row[0,29] = 20.75070391
$row[0,30] = 10.7549154
$row[0,31] = 0.85456375
$row[0,32] = 9.074045549999999
$row[0,33] = 9.856447360000001
$row[0,34] = 9.63845149
$row[0,35] = 38.79616001
$row[0,36] = 18.50886606
$row[0,37] = 8.48410778
$row[0,38] = 10.2199707
$row[0,39] = 37.39990613
$row[0,40] = 17.81906869
$row[0,41] = 5.01964729
$row[0,43] = 13.51031079
$row[0,44] = 26.74159645
$row[0,45] = 9.82216519
$row[0,46] = 20.25559872
$row[0,47] = 21.4266656
$row[0,48] = 7.75808352
$row[0,49] = 5.27235046
$row[0,51] = 10.83813674
$row[0,52] = 3.62592391
$row[0,53] = 31.8996932
$row[0,54] = 0
$row[0,55] = 0
$ws.Range("B142:BE142").Value = $row

# Row 143: fill data columns B:BE (columns E, N, AC, AR, AZ stay blank, as in every other data row)
$row = New-Object 'object[,]' 1,56
$row[0,0] = 12.44247602
$row[0,1] = 27.53835207
$row[0,2] = 24.70792822
$row[0,4] = 36.74927442
$row[0,5] = 7.87585451
$row[0,6] = 11.80438755
$row[0,7] = 8.5074387
$row[0,8] = 0
$row[0,9] = 12.32107749
$row[0,10] = 17.30358404
$row[0,11] = 18.57975647
$row[0,13] = 0
$row[0,14] = 48.60598331
$row[0,15] = 37.32640142
$row[0,16] = 9.61809775
$row[0,17] = 15.30158142
$row[0,18] = 34.05139093
$row[0,19] = 14.31646152
$row[0,20] = 15.35497133
$row[0,21] = 9.346437740000001
$row[0,22] = 4.68614514
$row[0,23] = 10.33372303
$row[0,24] = 6.83991943
$row[0,25] = 16.19301166
$row[0,26] = 26.58231066
$row[0,28] = 30.42306606
$row[0,29] = 15.32387825
$row[0,30] = 12.65147811
$row[0,31] = 0
$row[0,32] = 22.67297304
$row[0,33] = 27.42967304
$row[0,34] = 11.04227034
$row[0,35] = 30.15879603
$row[0,36] = 23.79927476
$row[0,37] = 7.52842621
$row[0,38] = 9.37713857
$row[0,39] = 38.32405913
$row[0,40] = 17.61368456
$row[0,41] = 3.06749394
$row[0,43] = 9.69094392
$row[0,44] = 37.10950704
$row[0,45] = 5.28807802
$row[0,46] = 18.31907826
$row[0,47] = 22.45793246
$row[0,48] = 4.87906805
$row[0,49] = 6.70060031
$row[0,51] = 8.4187335
$row[0,52] = 8.377298809999999
$row[0,53] = 32.75413663
$row[0,54] = 0
$row[0,55] = 0
$ws.Range("B143:BE143").Value = $row

# Row 144: fill data columns B:BE (columns E, N, AC, AR, AZ stay blank, as in every other data row)
$row = New-Object 'object[,]' 1,56
$row[0,0] = 8.30925045
$row[0,1] = 27.58211666
$row[0,2] = 35.30483846
$row[0,4] = 38.67392941
$row[0,5] = 8.651165799999999
$row[0,6] = 17.55269651
$row[0,7] = 5.69052992
$row[0,8] = 0
$row[0,9] = 8.95483595
$row[0,10] = 17.19627665
$row[0,11] = 21.32011093
$row[0,13] = 0
$row[0,14] = 38.83660193
$row[0,15] = 38.43167252
$row[0,16] = 10.31799186
$row[0,17] = 20.01096584
$row[0,18] = 39.83672276
$row[0,19] = 17.40496595
$row[0,20] = 26.92229405
$row[0,21] = 9.55353895
$row[0,22] = 6.13074275
$row[0,23] = 8.767332529999999
$row[0,24] = 6.1041072
$row[0,25] = 13.97322583
$row[0,26] = 28.18982487
$row[0,28] = 34.07656066
$row[0,29] = 10.44772007
$row[0,30] = 13.62443181
$row[0,31] = 0
$row[0,32] = 46.47372249
$row[0,33] = 21.97812319
$row[0,34] = 9.287312890000001
$row[0,35] = 22.11509164
$row[0,36] = 27.39474184
$row[0,37] = 6.46369913
$row[0,38] = 9.37971368
$row[0,39] = 41.6575904
$row[0,40] = 16.75869192
$row[0,41] = 3.45043751
$row[0,43] = 6.32914235
$row[0,44] = 37.0947152
$row[0,45] = 1.70276781
$row[0,46] = 16.05508527
$row[0,47] = 25.6062826
$row[0,48] = 10.69838389
$row[0,49] = 6.18870566
$row[0,51] = 6.08181765
$row[0,52] = 14.86868986
$row[0,53] = 34.12045205
$row[0,54] = 0
$row[0,55] = 0
$ws.Range("B144:BE144").Value = $row

# New trailing rows: date label only in column A, no numeric data yet
$ws.Range("A145").Value = "23 06 2020"
$ws.Range("A146").Value = "24 06 2020"
$ws.Range("A147").Value = "25 06 2020"
$ws.Range("A148").Value = "26 06 2020"

